# PNAD 2009 - correção nos dados e inicio da analise
# The sheet had three "header-only" rows (no numeric data) that are removed,
# shifting every data row up. The "urbana" row's values are also duplicated
# into "rural" (replacing rural's previous numbers), and the "unnamed:
# 5_level_1" column header becomes "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 2: "unnamed: 5_level_1" -> "total"
$ws.Range("F2").Value = "total"

# Remove the three label-only rows (no data), bottom-up so row numbers for
# rows not yet processed stay valid:
#   row 41 -> "fonte: ibge, ..." footnote row
#   row 8  -> "grandes regiões e unidades da federação" section header
#   row 5  -> "situação do domicílio" section header
$ws.Range("A41").EntireRow.Delete()
$ws.Range("A8").EntireRow.Delete()
$ws.Range("A5").EntireRow.Delete()

# After the deletions above, the old "rural" row (now row 6) still carries
# its own old figures; the corrected data duplicates "urbana"'s values here.
$ws.Range("B6").Value = 0.63
$ws.Range("C6").Value = 0.72
$ws.Range("D6").Value = 0.86
$ws.Range("E6").Value = 0.31
$ws.Range("F6").Value = 0.47
$ws.Range("G6").Value = 0.7
